$d = $word.ActiveDocument

$replacements = @(
    @("722÷7=", "813÷9="),
    @("878÷8=", "258÷6="),
    @("748÷9=", "406÷8="),
    @("133÷2=", "540÷5="),
    @("167÷9=", "973÷9="),
    @("781÷2=", "389÷7="),
    @("966÷2=", "509÷6="),
    @("419÷6=", "735÷2="),
    @("758÷2=", "440÷7="),
    @("819÷9=", "204÷8="),
    @("807÷5=", "101÷8="),
    @("870÷3=", "858÷3="),
    @("370÷3=", "839÷8="),
    @("229÷6=", "819÷2="),
    @("682÷7=", "185÷7="),
    @("992÷7=", "670÷4="),
    @("108÷6=", "282÷5="),
    @("852÷8=", "637÷5="),
    @("724÷8=", "694÷7="),
    @("706÷9=", "784÷9="),
    @("144÷2=", "924÷2="),
    @("396÷2=", "159÷7="),
    @("284÷7=", "206÷2="),
    @("997÷3=", "730÷4="),
    @("656÷3=", "163÷2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
